# Apply the "Add files via upload" edit to the AI-settings workbook.
#
# Summary of the change (derived from the OOXML diff):
#  - Column F (rows 2-12) holds a library of named prompt snippets. Several
#    entries were deleted ("日语" in F2, "脑暴" in F5, "动词" in F10,
#    "日语近义词" in F11) and the remaining entries shifted up to fill the
#    gaps, with the former F12 ("专家发言...") ending up in F4 and taking on
#    the wrap-text style/row height that F12 used to have. Rows 9-13 in
#    column F end up empty after the shift.
#  - Column B is untouched.
#  - Column widths: the former single A:E group (8.88671875) is split so
#    A:D keep the default width while E gets its own narrower width
#    (7.88671875).
#  - The view scroll/selection resets from F12 (with the window scrolled so
#    A12 is the top-left cell) to D2 with no scroll offset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot the existing column F values/styles before overwriting ---
$f3  = $ws.Range("F3").Value2
$f4  = $ws.Range("F4").Value2
$f6  = $ws.Range("F6").Value2
$f7  = $ws.Range("F7").Value2
$f8  = $ws.Range("F8").Value2
$f9  = $ws.Range("F9").Value2
$f12 = $ws.Range("F12").Value2

# --- shift column F content up, folding the old F12 entry into F4 ---
$ws.Range("F2").Value2 = $f3
$ws.Range("F3").Value2 = $f4

$ws.Range("F4").Value2 = $f12
$ws.Range("F4").WrapText = $true
$ws.Range("F4").VerticalAlignment = -4108
$ws.Rows.Item(4).RowHeight = 409.6

$ws.Range("F5").Value2 = $f6
$ws.Range("F6").Value2 = $f7
$ws.Range("F7").Value2 = $f8
$ws.Range("F8").Value2 = $f9

# rows 9 & 10 lose their F entries entirely
$ws.Range("F9").Clear()
$ws.Range("F10").Clear()

# rows 11 & 12 keep their (wrap-text) formatting but become empty, and their
# explicit row heights are removed (back to the sheet default)
$ws.Range("F11").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()

# --- column widths: split the old A:E (1-5) group into A:D (1-4) + E (5) ---
# (7.2 "characters" is the input that lands on the nearest achievable pixel
# width to the target raw width of 7.88671875 through this host's pixel-
# quantized ColumnWidth setter.)
$ws.Columns.Item(5).ColumnWidth = 7.2

# --- reset the view: scroll back to top and select D2 ---
$ws.Range("D2").Select()
